$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows cell G2 ("ticketDescriptionHighlighting" for the second
# ticket row) being re-saved: its text content is byte-for-byte identical
# before and after, but Excel re-serialized the shared-strings table and the
# entry for this particular string ended up re-appended at the end of the
# table (sharedStrings index 30 -> 32), which is why G2/G3/G4's <v> shared
# string indices all shift down by one. This happens when the cell is
# re-entered/re-committed in Excel without changing its visible text.
#
# Reproduce that user action: re-commit G2's existing value unchanged.
$g2 = $ws.Range("G2")
$originalText = $g2.Value2
$g2.Value = $originalText
